$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# B5: "Tóth Dániel" -> "Tóth Dániel, Jakab Botond"
$ws.Range("B5").Value = "Tóth Dániel, Jakab Botond"

# A8: "Vue Design elkezdése" -> "Vue navbar design"
$ws.Range("A8").Value = "Vue navbar design"

# A10: "Vue navbar design" -> "Vue navbar létrehoz"
$ws.Range("A10").Value = "Vue navbar létrehoz"

# B11: (empty) -> "Jakab Botond, Krausz Márton, Tóth Dániel"
$ws.Range("B11").Value = "Jakab Botond, Krausz Márton, Tóth Dániel"
